# Update "想去人数" (want-to-go count) figures in column F for the two
# sheets that carry the full event list ("展览" and "全部类型").
# Each entry maps a row number to its new value.

$updates = @{
    2  = 8845
    3  = 8247
    4  = 144
    5  = 197
    9  = 158
    12 = 748
    13 = 205
    14 = 5324
    17 = 3
    18 = 23
    21 = 154
    22 = 8
}

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
